# Add season-record columns (Wins / Losses / Ties) to the player table.
# (We downloaded only team statistics before; this fills in each player's
# team season record next to their stats.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns AD, AE, AF -- clone the existing header
# format (bold, centered, bordered) from AC1 so the new headers match the
# rest of row 1, then set the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-52: every player gets their team's season record (75-87-0).
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 75
    $ws.Cells.Item($row, 31).Value = 87
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-52"
